$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.141.47"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "'2.638.68"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'594.31"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("D6").Value = "'158.45"
$ws.Range("E6").Value = "  +1.85%  "

$ws.Range("D8").Value = "'0.542"
$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("D11").Value = "'5.26"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "'0.349"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").Value = "'27.85"
$ws.Range("E13").Value = "  -1.95%  "

$ws.Range("D14").Value = "'3.117.89"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").Value = "'0.0000187"
$ws.Range("E15").Value = "  -4.12%  "

$ws.Range("D16").Value = "'68.046.89"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "'2.641.65"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").Value = "'11.31"
$ws.Range("E18").Value = "  -2.36%  "

$ws.Range("D19").Value = "'358.83"
$ws.Range("E19").Value = "  -2.51%  "

$ws.Range("D20").Value = "'7.31"
$ws.Range("E20").Value = "  -3.62%  "

$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  -2.43%  "

$ws.Range("D22").Value = "'4.74"
$ws.Range("E22").Value = "  -3.80%  "

$ws.Range("D23").Value = "'2.06"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("D24").Value = "'74.66"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").Value = "'9.75"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").Value = "'2.773.01"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "'0.0000103"
$ws.Range("E28").Value = "  -5.58%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "'559.34"
$ws.Range("E30").Value = "  -4.19%  "

$ws.Range("D31").Value = "'7.94"
$ws.Range("E31").Value = "  -3.23%  "

$ws.Range("D32").Value = "'1.38"
$ws.Range("E32").Value = "  -4.59%  "

$ws.Range("D33").Value = "'1.84"
$ws.Range("E33").Value = "  -1.99%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'0.127"
$ws.Range("E35").Value = "  -4.48%  "

$ws.Range("D36").Value = "'1.54"
$ws.Range("E36").Value = "  -4.83%  "

$ws.Range("D37").Value = "'159.70"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").Value = "'19.66"
$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("D39").Value = "'0.368"
$ws.Range("E39").Value = "  -2.12%  "

$ws.Range("D40").Value = "'1.86"
$ws.Range("E40").Value = "  -2.86%  "

$ws.Range("D41").Value = "'5.27"
$ws.Range("E41").Value = "  -3.48%  "

$ws.Range("D42").Value = "'17.79"
$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("D43").Value = "'2.58"
$ws.Range("E43").Value = "  -5.58%  "

$ws.Range("D44").Value = "'0.0₆0325"
$ws.Range("E44").Value = "  -3.44%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "'156.83"
$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").Value = "'3.72"
$ws.Range("E47").Value = "  -2.59%  "

$ws.Range("D48").Value = "'21.82"
$ws.Range("E48").Value = "  -1.52%  "

$ws.Range("E49").Value = "  -3.41%  "

$ws.Range("D50").Value = "'0.0773"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").Value = "'0.611"
$ws.Range("E51").Value = "  -1.39%  "
